# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E35) previously listed the 20 monthly
# periods in descending order (1908 .. 1801). The database was updated so
# the periods now run in ascending chronological order (1801 .. 1908), and
# the partial-month "Valor Mora" amount (126667, vs. the standard 200000)
# now follows period 1908 instead of period 1801.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# The partial-period amount moves from the (old) first row to the (new) last row.
$ws.Range("F16").Value = 200000
$ws.Range("F35").Value = 126667
